# "Pruebas sobre guardado de archivos"
# The course year shown in cell A2 ("Curso" / A1 label) is updated from the
# upcoming year to the previous one, and the active selection on the sheet
# moves from C3 to A3 (matching the cursor position recorded in the saved
# file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held "2024-25" -> change it to "2023-24".
$ws.Range("A2").Value = "2023-24"

# Move/save the selection to A3 (was C3 before the edit).
[void]$ws.Range("A3").Select()
